# Natmi following Dr Hou advice:
# Recompute the Il6 -> Il6ra sending/target cluster table to include "ECs"
# as a sending cluster (in addition to the existing FAPs and sCs), giving
# the full 3x3 cross-product of sending/target clusters (ECs, FAPs, sCs)
# instead of the previous 2x3. Existing rows 2-7 get refreshed numbers and
# three new rows (8-10, sending cluster = sCs) are appended so the sheet
# now spans A1:T10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Il6", "Il6ra", "ECs", 2, 0.6666666666666666, 5.67446, 17.02338, 0.245923038656222, 0.245923038656222, 3, 1, 14.39487933333333, 43.184638, 0.6893176795959761, 0.6893176795959762, 81.68316698182667, 735.14850283644, 0.1695190983656985, 0.1695190983656985),
  @("ECs", "Il6", "Il6ra", "FAPs", 2, 0.6666666666666666, 5.67446, 17.02338, 0.245923038656222, 0.245923038656222, 3, 1, 5.526052333333333, 16.578157, 0.264622264871545, 0.2646222648715451, 31.35736292340667, 282.21626631066, 0.06507671147330199, 0.06507671147330202),
  @("ECs", "Il6", "Il6ra", "sCs", 2, 0.6666666666666666, 5.67446, 17.02338, 0.245923038656222, 0.245923038656222, 3, 1, 0.9618626666666668, 2.885588, 0.04606005553247879, 0.04606005553247879, 5.458051227493334, 49.12246104744001, 0.01132722881722151, 0.01132722881722152),
  @("FAPs", "Il6", "Il6ra", "ECs", 3, 1, 14.93150033333333, 44.794501, 0.647110021688359, 0.647110021688359, 3, 1, 14.39487933333333, 43.184638, 0.6893176795959761, 0.6893176795959762, 214.9371455639598, 1934.434310075638, 0.4460643785935214, 0.4460643785935214),
  @("FAPs", "Il6", "Il6ra", "FAPs", 3, 1, 14.93150033333333, 44.794501, 0.647110021688359, 0.647110021688359, 3, 1, 5.526052333333333, 16.578157, 0.264622264871545, 0.2646222648715451, 82.5122522571841, 742.610270314657, 0.1712397195602482, 0.1712397195602482),
  @("FAPs", "Il6", "Il6ra", "sCs", 3, 1, 14.93150033333333, 44.794501, 0.647110021688359, 0.647110021688359, 3, 1, 0.9618626666666668, 2.885588, 0.04606005553247879, 0.04606005553247879, 14.36205272795422, 129.258474551588, 0.02980592353458936, 0.02980592353458937),
  @("sCs", "Il6", "Il6ra", "ECs", 3, 1, 2.468169, 7.404507, 0.106966939655419, 0.106966939655419, 3, 1, 14.39487933333333, 43.184638, 0.6893176795959761, 0.6893176795959762, 35.528994929274, 319.760954363466, 0.07373420263675623, 0.07373420263675623),
  @("sCs", "Il6", "Il6ra", "FAPs", 3, 1, 2.468169, 7.404507, 0.106966939655419, 0.106966939655419, 3, 1, 5.526052333333333, 16.578157, 0.264622264871545, 0.2646222648715451, 13.639231061511, 122.753079553599, 0.02830583383799486, 0.02830583383799487),
  @("sCs", "Il6", "Il6ra", "sCs", 3, 1, 2.468169, 7.404507, 0.106966939655419, 0.106966939655419, 3, 1, 0.9618626666666668, 2.885588, 0.04606005553247879, 0.04606005553247879, 2.374039616124, 21.366356545116, 0.004926903180667907, 0.004926903180667907)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($startRow + $i, $c).Value = $row[$c-1]
  }
}

